$d = $word.ActiveDocument

# Remove the now-unused custom "OA_style_*" paragraph styles
# (Gold / Green / Gray) from the style sheet.
# Delete in reverse definition order (Gray, Green, Gold) so that each
# deletion doesn't invalidate the index used to look up the next style.
$stylesToRemove = @("OA_style_Gray", "OA_style_Green", "OA_style_Gold")

foreach ($styleName in $stylesToRemove) {
    $style = $d.Styles($styleName)
    if ($style -ne $null) {
        $style.Delete()
    }
}
